# This edit rotates a block of species-observation rows (rows 52-58) in the
# "Artfynd" worksheet. The data for rows 52, 53, 54, 56 and 58 gets replaced
# by the data that used to sit in the row above it in that same cyclic
# group (52<-58, 53<-52, 54<-53, 56<-54, 58<-56); rows 55 and 57 act as
# fixed points and keep their original data.
#
# Only the columns that actually differ between these rows need to be
# touched: A, B, D, E, F, G, H, P, Q, R and AC (all the other columns -
# C, I, J, K, N, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY -
# are identical across the whole block, so leaving them untouched keeps the
# row's effective content correct).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "before" values for every row involved, captured up front
# so that writes to one row don't clobber the source data needed for the
# next one.
$rows = @(52, 53, 54, 55, 56, 57, 58)
$snapshot = @{}

foreach ($r in $rows) {
    $snapshot[$r] = @{
        A  = $ws.Cells.Item($r, 1).Value2   # Id
        B  = $ws.Cells.Item($r, 2).Value2   # Taxonsorteringsordning
        D  = $ws.Cells.Item($r, 4).Value2   # Rödlistade
        E  = $ws.Cells.Item($r, 5).Value2   # TaxonId
        F  = $ws.Cells.Item($r, 6).Value2   # Artnamn
        G  = $ws.Cells.Item($r, 7).Value2   # Vetenskapligt namn
        H  = $ws.Cells.Item($r, 8).Value2   # Auktor
        P  = $ws.Cells.Item($r, 16).Value2  # Lokalnamn
        Q  = $ws.Cells.Item($r, 17).Value2  # Ost
        R  = $ws.Cells.Item($r, 18).Value2  # Nord
        AC = $ws.Cells.Item($r, 29).Value2  # Publik kommentar
    }
}

# Destination row -> source row (the row whose data it should now hold).
$mapping = @{
    52 = 58
    53 = 52
    54 = 53
    56 = 54
    58 = 56
}

foreach ($dest in @(52, 53, 54, 56, 58)) {
    $src = $mapping[$dest]
    $data = $snapshot[$src]

    $ws.Cells.Item($dest, 1).Value = $data.A
    $ws.Cells.Item($dest, 2).Value = $data.B
    $ws.Cells.Item($dest, 4).Value = $data.D
    $ws.Cells.Item($dest, 5).Value = $data.E
    $ws.Cells.Item($dest, 6).Value = $data.F
    $ws.Cells.Item($dest, 7).Value = $data.G
    $ws.Cells.Item($dest, 8).Value = $data.H
    $ws.Cells.Item($dest, 16).Value = $data.P
    $ws.Cells.Item($dest, 17).Value = $data.Q
    $ws.Cells.Item($dest, 18).Value = $data.R

    if ($data.AC) {
        $ws.Cells.Item($dest, 29).Value = $data.AC
    } else {
        $ws.Cells.Item($dest, 29).Value = ""
    }
}
